{"js": "// The document's Objetivos/Programa resumido/Programa/Avalia\u00e7\u00e3o/Bibliografia/Docente\n// sections get their body paragraphs' text content rotated between fixed slots\n// (paragraph styles, run formatting, and paragraph order all stay exactly the same;\n// only the w:t text that lives in each slot changes). We target each paragraph by\n// its fixed position in the body and overwrite its text in place, which preserves\n// every paragraph's style (Heading2/ListBullet/etc.) and run formatting (bold/italic).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Whole-paragraph (single run) text replacements, addressed by fixed paragraph index.\nconst paraReplacements = [\n  [5, \"Participa\u00e7\u00e3o em processo seletivo ou indica\u00e7\u00e3o de institui\u00e7\u00e3o para realiza\u00e7\u00e3o de est\u00e1gio. Submiss\u00e3o do plano de trabalho espec\u00edfico. Realiza\u00e7\u00e3o do est\u00e1gio e entrega do relat\u00f3rio de est\u00e1gio.\"],\n  [6, \"Participation in the selection process or indication of an institution to carry out an internship. Submission of the specific work plan. Conducting the internship and delivering the internship report.\"],\n  [11, \"Offer the opportunity to carry out professional training in a company or research institution, under the supervision of a professor from the Materials Engineering Department at EEL. Complement the general curricular training and psychologically and socially adapt the student to his/her future professional activity.\"],\n  [10, \"Participa\u00e7\u00e3o do aluno em processo seletivo de empresas, institui\u00e7\u00f5es de pesquisa ou no setor acad\u00eamico. O est\u00e1gio ser\u00e1 realizado sob a supervis\u00e3o de docente designado pela Comiss\u00e3o de Curso de Engenharia F\u00edsica. O conte\u00fado ser\u00e1 estabelecido no Plano de Trabalho entre o supervisor respons\u00e1vel pelo Est\u00e1gio e o docente supervisor. Apresenta\u00e7\u00e3o de relat\u00f3rio final sobre as atividades desenvolvidas no est\u00e1gio.\"],\n  [13, \"Supervis\u00e3o das atividades desenvolvidas pelo aluno durante o est\u00e1gio.\"],\n  [18, \"1176388 - Luiz Tadeu Fernandes Eleno\"],\n  [8, \"Oferecer oportunidade de realiza\u00e7\u00e3o de treinamento profissional em empresa ou institui\u00e7\u00e3o de pesquisa, sob supervis\u00e3o de docente do Departamento de Engenharia de Materiais da EEL. Complementar a forma\u00e7\u00e3o geral curricular e adaptar psicol\u00f3gica e socialmente o estudante \u00e0 sua futura atividade profissional.\"],\n];\n\nfor (const [index, newText] of paraReplacements) {\n  paragraphs.items[index].insertText(newText, \"Replace\");\n}\nawait context.sync();\n\n// Paragraph 16 (\"Avalia\u00e7\u00e3o\" bullet) holds three separate label/value run pairs\n// (M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o) inside one paragraph. Only the value\n// runs move; the bold labels and the w:br line breaks must stay untouched, so we\n// search-and-replace scoped to that paragraph's range instead of overwriting the\n// whole paragraph. Processed last-value-first so each search target is still\n// unique when it is looked up.\nconst p16 = paragraphs.items[16].getRange();\n\nasync function replaceWithin(range, oldText, newText) {\n  const results = range.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nawait replaceWithin(p16, \"Devido \u00e0s caracter\u00edsticas da disciplina, n\u00e3o ser\u00e1 oferecida recupera\u00e7\u00e3o.\", \"A ser definida com o orientador em fun\u00e7\u00e3o das atividades desenvolvidas no est\u00e1gio.\");\nawait replaceWithin(p16, \"A nota final ser\u00e1 baseada em relat\u00f3rio final e no desempenho no est\u00e1gio, a ser atribu\u00edda pelo docente supervisor do est\u00e1gio.\", \"Devido \u00e0s caracter\u00edsticas da disciplina, n\u00e3o ser\u00e1 oferecida recupera\u00e7\u00e3o.\");\nawait replaceWithin(p16, \"Supervis\u00e3o das atividades desenvolvidas pelo aluno durante o est\u00e1gio.\", \"A nota final ser\u00e1 baseada em relat\u00f3rio final e no desempenho no est\u00e1gio, a ser atribu\u00edda pelo docente supervisor do est\u00e1gio.\");\n", "ps1": "# The document's Objetivos/Programa resumido/Programa/Avalia\u00e7\u00e3o/Bibliografia/Docente\n# sections get their paragraphs' text content rotated between fixed slots (paragraph\n# styles, run formatting, and paragraph order all stay exactly the same; only the\n# text that lives in each slot changes). We target each paragraph by its fixed\n# position in the document and overwrite its Range.Text in place, which preserves\n# that paragraph's style (Heading2/ListBullet/etc.) and run formatting.\n\n$d = $word.ActiveDocument\n\n# Whole-paragraph (single run) text replacements, addressed by fixed 1-based\n# paragraph index.\n$d.Paragraphs.Item(6).Range.Text  = \"Participa\u00e7\u00e3o em processo seletivo ou indica\u00e7\u00e3o de institui\u00e7\u00e3o para realiza\u00e7\u00e3o de est\u00e1gio. Submiss\u00e3o do plano de trabalho espec\u00edfico. Realiza\u00e7\u00e3o do est\u00e1gio e entrega do relat\u00f3rio de est\u00e1gio.\"\n$d.Paragraphs.Item(7).Range.Text  = \"Participation in the selection process or indication of an institution to carry out an internship. Submission of the specific work plan. Conducting the internship and delivering the internship report.\"\n$d.Paragraphs.Item(12).Range.Text = \"Offer the opportunity to carry out professional training in a company or research institution, under the supervision of a professor from the Materials Engineering Department at EEL. Complement the general curricular training and psychologically and socially adapt the student to his/her future professional activity.\"\n$d.Paragraphs.Item(11).Range.Text = \"Participa\u00e7\u00e3o do aluno em processo seletivo de empresas, institui\u00e7\u00f5es de pesquisa ou no setor acad\u00eamico. O est\u00e1gio ser\u00e1 realizado sob a supervis\u00e3o de docente designado pela Comiss\u00e3o de Curso de Engenharia F\u00edsica. O conte\u00fado ser\u00e1 estabelecido no Plano de Trabalho entre o supervisor respons\u00e1vel pelo Est\u00e1gio e o docente supervisor. Apresenta\u00e7\u00e3o de relat\u00f3rio final sobre as atividades desenvolvidas no est\u00e1gio.\"\n$d.Paragraphs.Item(14).Range.Text = \"Supervis\u00e3o das atividades desenvolvidas pelo aluno durante o est\u00e1gio.\"\n$d.Paragraphs.Item(19).Range.Text = \"1176388 - Luiz Tadeu Fernandes Eleno\"\n$d.Paragraphs.Item(9).Range.Text  = \"Oferecer oportunidade de realiza\u00e7\u00e3o de treinamento profissional em empresa ou institui\u00e7\u00e3o de pesquisa, sob supervis\u00e3o de docente do Departamento de Engenharia de Materiais da EEL. Complementar a forma\u00e7\u00e3o geral curricular e adaptar psicol\u00f3gica e socialmente o estudante \u00e0 sua futura atividade profissional.\"\n\n# Paragraph 17 (\"Avalia\u00e7\u00e3o\" bullet) holds three separate label/value run pairs\n# (M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o) inside one paragraph. Only the value\n# runs move; the bold labels and the line breaks between them must stay untouched,\n# so we use Find/Replace scoped to that paragraph's range instead of overwriting\n# the whole paragraph. A fresh Range is re-fetched before each call because\n# Find.Execute collapses its range down to the replacement text. Processed\n# last-value-first so each search target is still unique in the document when\n# it is looked up.\n$d.Paragraphs.Item(17).Range.Find.Execute(\"Devido \u00e0s caracter\u00edsticas da disciplina, n\u00e3o ser\u00e1 oferecida recupera\u00e7\u00e3o.\", $false, $true, $false, $false, $false, $true, 1, $false, \"A ser definida com o orientador em fun\u00e7\u00e3o das atividades desenvolvidas no est\u00e1gio.\", 2)\n$d.Paragraphs.Item(17).Range.Find.Execute(\"A nota final ser\u00e1 baseada em relat\u00f3rio final e no desempenho no est\u00e1gio, a ser atribu\u00edda pelo docente supervisor do est\u00e1gio.\", $false, $true, $false, $false, $false, $true, 1, $false, \"Devido \u00e0s caracter\u00edsticas da disciplina, n\u00e3o ser\u00e1 oferecida recupera\u00e7\u00e3o.\", 2)\n$d.Paragraphs.Item(17).Range.Find.Execute(\"Supervis\u00e3o das atividades desenvolvidas pelo aluno durante o est\u00e1gio.\", $false, $true, $false, $false, $false, $true, 1, $false, \"A nota final ser\u00e1 baseada em relat\u00f3rio final e no desempenho no est\u00e1gio, a ser atribu\u00edda pelo docente supervisor do est\u00e1gio.\", 2)\n"}
